$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 403; existing rows 403:423 shift down to 407:427
$ws.Rows("403:406").Insert()

# New weekly data block (week of 2022-05-25) for Pera at Vega Monumental Concepcion
$fecha = Get-Date -Year 2022 -Month 5 -Day 25 -Hour 0 -Minute 0 -Second 0
$newRows = @(
    @{ Row=403; Fecha=$fecha; Variedad="Abate Fettel";       Calidad="Primera"; Vol=50; Min=8000; Max=8000; Prom=8000; Precio=500 },
    @{ Row=404; Fecha=$fecha; Variedad="Abate Fettel";       Calidad="Segunda"; Vol=50; Min=7000; Max=7000; Prom=7000; Precio=438 },
    @{ Row=405; Fecha=$fecha; Variedad="Packham's Triumph";  Calidad="Primera"; Vol=50; Min=8000; Max=8000; Prom=8000; Precio=500 },
    @{ Row=406; Fecha=$fecha; Variedad="Packham's Triumph";  Calidad="Segunda"; Vol=50; Min=7000; Max=7000; Prom=7000; Precio=438 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 11
    $ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($row, 3).Value = "Bíobío"
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = 8
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100104
    $ws.Cells.Item($row, 8).Value = "Frutos de pepita"
    $ws.Cells.Item($row, 9).Value = 100104005
    $ws.Cells.Item($row, 10).Value = "Pera"
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Vol
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = "`$/caja 16 kilos empedrada"
    $ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 19).Value = $r.Precio
    $ws.Cells.Item($row, 20).Value = 16
}
